$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C1").Value = "Dept"
$ws.Range("C2").Value = "AP"
[void]$ws.Range("D10").Select()
